# 26-10-2023_10-00 AM.xlsx -- "subject and year strength update"
#
# 1) Remove the worksheet "G13" entirely (sheet + its shared strings
#    usages go away with it).
# 2) Rename the subject labels used in row 4 of every remaining sheet:
#       "D4CSE_CC"  -> "D3CSE_URN"
#       "D4CSE_NSC" -> "D4CSE_DWDM"
# 3) Update the year-wise strength figures (columns B, C, E, F) on the
#    remaining sheets: the "22xx0xx"-style numbers shrink from
#    2222xxx/3333xxx down to 21xxx/1111xxx (same trailing 3 digits).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1) Delete sheet "G13"
# ---------------------------------------------------------------------
$wb.Worksheets.Item("G13").Delete()

# ---------------------------------------------------------------------
# 2) Update the subject-label text on every remaining sheet (row 4)
# ---------------------------------------------------------------------
foreach ($name in @("S216", "G14", "G15", "S219")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("B4").Value = "D3CSE_URN"
    $ws.Range("E4").Value = "D3CSE_URN"

    $ws.Range("C4").Value = "D4CSE_DWDM"
    $ws.Range("F4").Value = "D4CSE_DWDM"
}

# ---------------------------------------------------------------------
# 3) Update the numeric strength figures (B/C/E/F columns)
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("S216")
$ws.Range("B5").Value = 21000
$ws.Range("C5").Value = 1111000
$ws.Range("E5").Value = 21011
$ws.Range("F5").Value = 1111009
$ws.Range("B6").Value = 21001
$ws.Range("C6").Value = 1111001
$ws.Range("E6").Value = 21012
$ws.Range("F6").Value = 1111010
$ws.Range("B7").Value = 21002
$ws.Range("C7").Value = 1111002
$ws.Range("E7").Value = 21013
$ws.Range("F7").Value = 1111011
$ws.Range("B8").Value = 21003
$ws.Range("C8").Value = 1111003
$ws.Range("E8").Value = 21014
$ws.Range("F8").Value = 1111012
$ws.Range("B9").Value = 21004
$ws.Range("C9").Value = 1111004
$ws.Range("E9").Value = 21015
$ws.Range("F9").Value = 1111013
$ws.Range("B10").Value = 21005
$ws.Range("C10").Value = 1111005
$ws.Range("E10").Value = 21016
$ws.Range("F10").Value = 1111014
$ws.Range("B11").Value = 21006
$ws.Range("C11").Value = 1111006
$ws.Range("E11").Value = 21017
$ws.Range("F11").Value = 1111015
$ws.Range("B12").Value = 21007
$ws.Range("C12").Value = 1111007
$ws.Range("E12").Value = 21018
$ws.Range("F12").Value = 1111016
$ws.Range("B13").Value = 21008
$ws.Range("C13").Value = 1111008
$ws.Range("E13").Value = 21019
$ws.Range("F13").Value = 1111017
$ws.Range("B14").Value = 21009
$ws.Range("B15").Value = 21010

$ws = $wb.Worksheets.Item("G14")
$ws.Range("B5").Value = 21020
$ws.Range("C5").Value = 1111018
$ws.Range("E5").Value = 21028
$ws.Range("F5").Value = 1111026
$ws.Range("B6").Value = 21021
$ws.Range("C6").Value = 1111019
$ws.Range("E6").Value = 21029
$ws.Range("F6").Value = 1111027
$ws.Range("B7").Value = 21022
$ws.Range("C7").Value = 1111020
$ws.Range("E7").Value = 21030
$ws.Range("F7").Value = 1111028
$ws.Range("B8").Value = 21023
$ws.Range("C8").Value = 1111021
$ws.Range("E8").Value = 21031
$ws.Range("F8").Value = 1111029
$ws.Range("B9").Value = 21024
$ws.Range("C9").Value = 1111022
$ws.Range("E9").Value = 21032
$ws.Range("F9").Value = 1111030
$ws.Range("B10").Value = 21025
$ws.Range("C10").Value = 1111023
$ws.Range("E10").Value = 21033
$ws.Range("F10").Value = 1111031
$ws.Range("B11").Value = 21026
$ws.Range("C11").Value = 1111024
$ws.Range("E11").Value = 21034
$ws.Range("F11").Value = 1111032
$ws.Range("B12").Value = 21027
$ws.Range("C12").Value = 1111025
$ws.Range("E12").Value = 21035
$ws.Range("F12").Value = 1111033

$ws = $wb.Worksheets.Item("G15")
$ws.Range("B5").Value = 21036
$ws.Range("C5").Value = 1111034
$ws.Range("E5").Value = 21043
$ws.Range("F5").Value = 1111041
$ws.Range("B6").Value = 21037
$ws.Range("C6").Value = 1111035
$ws.Range("E6").Value = 21044
$ws.Range("F6").Value = 1111042
$ws.Range("B7").Value = 21038
$ws.Range("C7").Value = 1111036
$ws.Range("E7").Value = 21045
$ws.Range("F7").Value = 1111043
$ws.Range("B8").Value = 21039
$ws.Range("C8").Value = 1111037
$ws.Range("E8").Value = 21046
$ws.Range("F8").Value = 1111044
$ws.Range("B9").Value = 21040
$ws.Range("C9").Value = 1111038
$ws.Range("E9").Value = 21047
$ws.Range("F9").Value = 1111045
$ws.Range("B10").Value = 21041
$ws.Range("C10").Value = 1111039
$ws.Range("E10").Value = 21048
$ws.Range("F10").Value = 1111046
$ws.Range("B11").Value = 21042
$ws.Range("C11").Value = 1111040
$ws.Range("E11").Value = 21049
$ws.Range("F11").Value = 1111047
$ws.Range("F12").Value = 1111048
$ws.Range("F13").Value = 1111049

$ws = $wb.Worksheets.Item("S219")
$ws.Range("B5").Value = 21050
$ws.Range("C5").Value = 1111050
$ws.Range("E5").Value = 21055
$ws.Range("F5").Value = 1111055
$ws.Range("B6").Value = 21051
$ws.Range("C6").Value = 1111051
$ws.Range("E6").Value = 21056
$ws.Range("F6").Value = 1111056
$ws.Range("B7").Value = 21052
$ws.Range("C7").Value = 1111052
$ws.Range("E7").Value = 21057
$ws.Range("F7").Value = 1111057
$ws.Range("B8").Value = 21053
$ws.Range("C8").Value = 1111053
$ws.Range("E8").Value = 21058
$ws.Range("F8").Value = 1111058
$ws.Range("B9").Value = 21054
$ws.Range("C9").Value = 1111054
$ws.Range("E9").Value = 21059
$ws.Range("F9").Value = 1111059
$ws.Range("E10").Value = 21060
$ws.Range("F10").Value = 1111060
$ws.Range("E11").Value = 21061
$ws.Range("F11").Value = 1111061
$ws.Range("E12").Value = 21062
$ws.Range("F12").Value = 1111062
